$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.35
$ws.Range("I2").Value = 2.88

$ws.Range("G10").Value = 1.85
$ws.Range("H10").Value = 3.35
$ws.Range("I10").Value = 4.15
$ws.Range("P10").Value = 1.45
$ws.Range("Q10").Value = 2.57
$ws.Range("T10").Value = 6
$ws.Range("U10").Value = 7.9
$ws.Range("W10").Value = 15
$ws.Range("X10").Value = 16
$ws.Range("AA10").Value = 6.5
$ws.Range("AE10").Value = 10.25
$ws.Range("AF10").Value = 22
$ws.Range("AI10").Value = 45

$ws.Range("G12").Value = 3.6
$ws.Range("H12").Value = 3.75
$ws.Range("I12").Value = 1.91
$ws.Range("K12").Value = 9.5
$ws.Range("U12").Value = 19
$ws.Range("V12").Value = 13
$ws.Range("X12").Value = 34
$ws.Range("AC12").Value = 67
$ws.Range("AF12").Value = 8.5

$ws.Range("G14").Value = 2.15
$ws.Range("I14").Value = 3.1
$ws.Range("N14").Value = 1.75
$ws.Range("O14").Value = 2.05
$ws.Range("T14").Value = 9.5
$ws.Range("U14").Value = 12
$ws.Range("AE14").Value = 11
$ws.Range("AG14").Value = 11

$ws.Range("G18").Value = 5.8
$ws.Range("H18").Value = 3.85
$ws.Range("I18").Value = 1.52
$ws.Range("L18").Value = 1.24
$ws.Range("M18").Value = 3.3
$ws.Range("N18").Value = 1.7
$ws.Range("O18").Value = 1.9
$ws.Range("S18").Value = 1.82
$ws.Range("T18").Value = 16.5
$ws.Range("U18").Value = 37
$ws.Range("V18").Value = 18
$ws.Range("W18").Value = 120
$ws.Range("X18").Value = 60
$ws.Range("Z18").Value = 11.25
$ws.Range("AA18").Value = 7.6
$ws.Range("AB18").Value = 16.5
$ws.Range("AF18").Value = 7.2
$ws.Range("AH18").Value = 10.75
$ws.Range("AI18").Value = 12

$ws.Range("G20").Value = 2.3
$ws.Range("H20").Value = 3.05
$ws.Range("I20").Value = 3.05
$ws.Range("M20").Value = 2.57
$ws.Range("N20").Value = 2.12
$ws.Range("T20").Value = 6.7
$ws.Range("U20").Value = 10.25
$ws.Range("V20").Value = 9.25
$ws.Range("W20").Value = 23
$ws.Range("X20").Value = 21
$ws.Range("AE20").Value = 8
$ws.Range("AF20").Value = 15
$ws.Range("AG20").Value = 11.25
$ws.Range("AH20").Value = 40
$ws.Range("AI20").Value = 30
$ws.Range("AJ20").Value = 40

$ws.Range("J24").Value = 1.06
$ws.Range("K24").Value = 10
$ws.Range("L24").Value = 1.33
$ws.Range("M24").Value = 3.25
$ws.Range("N24").Value = 2.08
$ws.Range("O24").Value = 1.73

$ws.Range("G34").Value = 1.36
$ws.Range("I34").Value = 6
$ws.Range("R34").Value = 1.67
$ws.Range("S34").Value = 2.1
$ws.Range("T34").Value = 11
$ws.Range("W34").Value = 10
$ws.Range("AC34").Value = 41
$ws.Range("AG34").Value = 19

$ws.Range("J40").Value = 1.04
$ws.Range("K40").Value = 13
$ws.Range("L40").Value = 1.22
$ws.Range("M40").Value = 4
$ws.Range("N40").Value = 1.7
$ws.Range("O40").Value = 2.1

$ws.Range("H45").Value = 4
$ws.Range("I45").Value = 1.8
$ws.Range("R45").Value = 1.67
$ws.Range("S45").Value = 2.1
$ws.Range("Y45").Value = 34
$ws.Range("AA45").Value = 7.5
$ws.Range("AE45").Value = 8.5

$ws.Range("H46").Value = 6.25
$ws.Range("AG46").Value = 11
